# timesheet.xlsx edit script
# - Update the book's tab-bar ratio (best effort; may not be honored by the
#   headless host, but set it anyway so the object model is consistent).
# - Move the active selection on Sheet1 from C7 to A9.
# - Change the "Total Time" formulas in column D from straight ABS(...)*24
#   (or the broken single-arg ROUND(...) in D3) into ROUND(ABS(...)*24, 1)
#   for every data row (3-14).
# - Populate the previously-empty timesheet entries for row 7 (11/15/2014)
#   and row 8 (11/16/2014) with a date, start time and end time, matching
#   the formatting already used by the rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window: widen the sheet-tab area ---
$win = $wb.Windows.Item(1)
$win.TabRatio = 0.423

# Reuse the exact same number-format codes already used by the sheet
# (note the escaped space in the time format, matching the existing
# styles.xml entry) so these new cells share the existing style indices
# instead of creating duplicate numFmt/cellXfs entries.
$dateFormat = "MM/DD/YY"
$timeFormat = 'HH:MM:SS\ AM/PM'

# --- New data for row 7 (date 11/15/2014, 23:30 -> 23:59:59) ---
$ws.Range("A7").Value = 41958
$ws.Range("A7").NumberFormat = $dateFormat
$ws.Range("B7").Value = 0.979166666666667
$ws.Range("B7").NumberFormat = $timeFormat
$ws.Range("C7").Value = 0.999988425925926
$ws.Range("C7").NumberFormat = $timeFormat

# --- New data for row 8 (date 11/16/2014, 01:30 -> 06:00) ---
$ws.Range("A8").Value = 41959
$ws.Range("A8").NumberFormat = $dateFormat
$ws.Range("B8").Value = 0.0625
$ws.Range("B8").NumberFormat = $timeFormat
$ws.Range("C8").Value = 0.25
$ws.Range("C8").NumberFormat = $timeFormat

# --- Fix up every "Total Time" formula in column D to round to 1 decimal ---
for ($row = 3; $row -le 14; $row++) {
    $cell = "D$row"
    $ws.Range($cell).Formula = "=ROUND(ABS(C$row-B$row) * 24, 1)"
}

# --- Move the active selection to A9 ---
$ws.Range("A9").Select() | Out-Null
